$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# Insert a new row at 13, shifting rows 13..23 down to 14..24.
# This matches the target layout where "Programa resumido:" (and
# everything below it) moves down one row, and the freed row 13
# becomes the new home for the "Docentes responsaveis" value.
$ws.Rows("13:13").Insert()

# Row 10 (Objetivos): replace the misplaced professor name with the
# real course-objectives paragraph.
$objText = @"
Visão integrada sobre petróleo e gás natural, desde a origem até o processamento primário. Descrições, características e aplicações dos derivados do petróleo. Processo e esquemas de refino e processamento do gás natural.
"@
$ws.Range("B10").Value = $objText
$ws.Range("C10").Value = $objText

# Row 13 (no label, sits under "Docentes responsaveis:"): the
# professor's name belongs here now.
$docente = "1285870 - Marcos Villela Barcza"
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente

# Row 14 (Programa resumido:): the long numbered syllabus outline.
$programaResumido = @"
1.Petróleo: histórico, constituinte, composição e classificação 
2.Geologia do petróleo: origem.
3.Prospecção de petróleo: métodos geológicos, potenciais, sísmicos; 
4.Perfuração: equipamentos, operações;
5.Completação e reservatórios: tipos, etapas, equipamentos, reservatórios;
6.Elevação: elevação natural, sistemas de bombeamentos;
7.Processamento primário: separação do gás natural, tratamento do óleo, tratamento da água, unidade de processamento de gás natural;
8.Derivados do petróleo: tipos, características, gás liquefeito de petróleo, gasolina automotiva, querosene de aviação, óleo diesel, óleos combustíveis industriais, óleos combustíveis marítimos, produtos especiais;
9.Processos de refino: objetivo, tipos de processos, esquemas de refino.
9.1- Destilação de petróleo: equipamentos, esquemas típicos, descrição e variáveis do processo;
9.2- Desasfaltação: carga, descrição e variáveis do processo, produtos;
9.3- Coqueamento retardado: carga, descrição e variáveis do processo, produtos;
9.4- Craqueamento catalítico: carga, descrição e variáveis do processo, produtos;
9.5- Hidrorrefino: carga, descrição e variáveis do processo, produtos;
9.6- Reforma catalítica: carga, descrição e variáveis do processo, produtos;
9.7- Alquilação e isomerização: carga, descrição e variáveis do processo, produtos;
9.8- Tratamento de derivados: tratamento com aminas, tratamentos cáusticos;
9.9- Geração de hidrogênio: carga, descrição e variáveis do processo;
9.10- Recuperação de Enxofre: Processo Claus.
10.Óleos básicos lubrificantes e parafinas: carga, descrição e variáveis do processo.
"@
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# Row 16 (Programa:): the single-paragraph syllabus text.
$programa = @"
Petróleo: histórico, constituinte, composição e classificação; Geologia do petróleo: origem; Prospecção de petróleo: métodos geológicos, potenciais, sísmicos; Perfuração: equipamentos, operações; Completação e reservatórios: tipos, etapas, equipamentos, reservatórios; Elevação: elevação natural, bombeios; Processamento primário: separação do gás natural, tratamento do óleo, tratamento da água, unidade de processamento de gás natural; Derivados do petróleo: tipos, características, gás liquefeito de petróleo, gasolina automotiva, querosene de aviação, óleo diesel, óleos combustíveis industriais, óleos combustíveis marítimos, produtos especiais; Processos de refino: objetivo, tipos de processos, esquemas de refino; Destilação de petróleo: equipamentos, esquemas típicos, descrição e variáveis do processo; Desasfaltação: carga, descrição e variáveis do processo, produtos; Coqueamento retardado: carga, descrição e variáveis do processo, produtos; Craqueamento catalítico: carga, descrição e variáveis do processo, produtos; Hidrorrefino: carga, descrição e variáveis do processo, produtos; Reforma catalítica: carga, descrição e variáveis do processo, produtos; Alquilação e isomerização: carga, descrição e variáveis do processo, produtos; Tratamento de derivados: tratamento com aminas, tratamentos cáusticos; Geração de hidrogênio: carga, descrição e variáveis do processo; Recuperação de Enxofre: Processo Claus; Óleos básicos lubrificantes e parafinas: carga, descrição e variáveis do processo.
"@
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Row 19 (Metodo:): teaching method description.
$metodo = @"
Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos e seminários
"@
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Row 20 (Criterio:): grading criteria description.
$criterio = @"
Provas, avaliação através de exercícios ou casos práticos elaborados fora de sala de aula.
"@
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Row 21 (Norma de recuperacao:): recovery-exam rule.
$norma = @"
Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação.
"@
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# Row 22 (Bibliografia:): reference list.
$biblio = @"
a)Speight, J. G., The Chemistry and Technology of Petroleum, CRC Press, 4ª Edição, 2007;
b)Thomas, J. E. (Organizador), Fundamentos de Engenharia de Petróleo, Editora Interciência, 2ª Edição, 2004;
c)Brasil, N. I., Araújo, M. A. S., Souza, E. C. M, Processamento de Petróleo e Gás, Editora LTC, 1ª Edição, 2012;
d)Fundamentos do Refino do Petróleo  Tecnologia e Economia, Szklo, A. S., Uller, V. C., Bonfá, M. H. P., Editora Interciência, 3ª Edição, 2012.
e)Oil and Gas Journal;
f)Revista Petro & Química.
"@
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio

# Column A's width definition originally (incorrectly) spanned
# columns 1-2; narrow it back to just column A.
$ws.Columns("A:A").ColumnWidth = 30.7109375

$true
